# Apply updated forecast / seasonality figures to the "Forecast Comparison"
# sheet and the roll-up totals on the "Summary" sheet.

$wb = $excel.ActiveWorkbook

# --- Forecast Comparison sheet -------------------------------------------
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Seasonality Index (column L)
$ws.Range("L2").Value  = 0.95
$ws.Range("L3").Value  = 1.14
$ws.Range("L4").Value  = 0.96
$ws.Range("L5").Value  = 0.8100000000000001
$ws.Range("L6").Value  = 1.05
$ws.Range("L7").Value  = 0.89
$ws.Range("L8").Value  = 0.83
$ws.Range("L9").Value  = 0.98
$ws.Range("L10").Value = 0.95
$ws.Range("L11").Value = 0.95
$ws.Range("L12").Value = 1.07
$ws.Range("L13").Value = 1.03
$ws.Range("L14").Value = 1.13
$ws.Range("L15").Value = 1.18
$ws.Range("L16").Value = 0.84
$ws.Range("L17").Value = 0.95

# MyForecast (column D)
$ws.Range("D7").Value  = 94
$ws.Range("D8").Value  = 93
$ws.Range("D10").Value = 91
$ws.Range("D12").Value = 85
$ws.Range("D13").Value = 80
$ws.Range("D14").Value = 69
$ws.Range("D15").Value = 68
$ws.Range("D16").Value = 62
$ws.Range("D17").Value = 62

# --- Summary sheet ---------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Summary")

# B9/B10 hold numeric-looking text values, so force text formatting before
# assigning the strings to avoid Excel auto-converting them to numbers.
$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "1277"

$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "672"
